# Day 24.1: manual scale override + geometry fallback + UI updates
#
# - Rates sheet gets explicit currency + overhead/contingency/profit/tax
#   knobs (INR / 5 / 3 / 10 / 18) instead of being blank.
# - Details sheet's INDIA row picks up a geometry-fallback estimate for
#   Materials/Labor (previously 0).
# - Summary / Charts / Compare sheets are recalculated downstream from
#   those inputs (materials/labor totals, cascading overhead ->
#   contingency -> profit -> tax -> grand total).
# - A couple of column widths are widened so the new (wider) INR figures
#   aren't clipped.

$wb = $excel.ActiveWorkbook

# Excel's ColumnWidth property (characters, Calibri 11 default font) is
# offset from the raw OOXML <col width> by 5/6 of a character; subtract
# that so the saved XML width lands on the exact target value.
$colWidthFudge = 5 / 6

# ---------------------------------------------------------------------
# Rates sheet: fill in the pricing knobs
# ---------------------------------------------------------------------
$wsRates = $wb.Worksheets.Item("Rates")
$wsRates.Cells.Item(4, 2).Value = "INR"
$wsRates.Cells.Item(5, 2).Value = 5
$wsRates.Cells.Item(6, 2).Value = 3
$wsRates.Cells.Item(7, 2).Value = 10
$wsRates.Cells.Item(8, 2).Value = 18

# ---------------------------------------------------------------------
# Details sheet: INDIA row geometry-fallback estimate + wider columns
# ---------------------------------------------------------------------
$wsDetails = $wb.Worksheets.Item("Details")
$wsDetails.Columns.Item(2).ColumnWidth = 19 - $colWidthFudge
$wsDetails.Columns.Item(3).ColumnWidth = 18 - $colWidthFudge
$wsDetails.Cells.Item(4, 2).Value = 342198.7871837872
$wsDetails.Cells.Item(4, 3).Value = 58990.1994645264

# ---------------------------------------------------------------------
# Summary sheet: currency label, wider amount column, recomputed totals
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Columns.Item(2).ColumnWidth = 20 - $colWidthFudge
$wsSummary.Range("A2").Value = "Currency: INR"
$wsSummary.Cells.Item(5, 2).Value = 345217.5371837872
$wsSummary.Cells.Item(6, 2).Value = 60436.6994645264
$wsSummary.Cells.Item(7, 2).Value = 405654.2366483136
$wsSummary.Cells.Item(8, 2).Value = 20282.71183241568
$wsSummary.Cells.Item(9, 2).Value = 12778.10845442188
$wsSummary.Cells.Item(10, 2).Value = 43871.50569351512
$wsSummary.Cells.Item(11, 2).Value = 86865.58127315991
$wsSummary.Cells.Item(12, 2).Value = 569452.1439018261

# ---------------------------------------------------------------------
# Charts sheet: Grand Total feeding the chart
# ---------------------------------------------------------------------
$wsCharts = $wb.Worksheets.Item("Charts")
$wsCharts.Cells.Item(3, 2).Value = 569452.1439018261

# ---------------------------------------------------------------------
# Compare sheet: INDIA column picks up the same geometry-fallback values
# ---------------------------------------------------------------------
$wsCompare = $wb.Worksheets.Item("Compare")
$wsCompare.Cells.Item(4, 2).Value = 342198.7871837872
$wsCompare.Cells.Item(5, 2).Value = 58990.1994645264
$wsCompare.Cells.Item(6, 2).Value = 401188.9866483136
$wsCompare.Cells.Item(10, 2).Value = 342198.7871837872
$wsCompare.Cells.Item(11, 2).Value = 58990.1994645264
